$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.421.25"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").Value = "1.965.65"
$ws.Range("E3").Value = "  -4.99%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'244.41"
$ws.Range("E5").Value = "  -3.30%  "
$ws.Range("D6").Value = "'0.618"
$ws.Range("E6").Value = "  -5.16%  "
$ws.Range("D7").Value = "'58.34"
$ws.Range("E7").Value = "  -10.94%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.371"
$ws.Range("E9").Value = "  -6.92%  "
$ws.Range("D10").Value = "'55.82"
$ws.Range("E10").Value = "  -6.18%  "
$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "'0.103"
$ws.Range("E12").Value = "  -1.25%  "
$ws.Range("D13").Value = "'22.19"
$ws.Range("E13").Value = "  -6.95%  "
$ws.Range("D14").Value = "'0.833"
$ws.Range("E14").Value = "  -10.04%  "
$ws.Range("D15").Value = "2.251.11"
$ws.Range("E15").Value = "  -5.17%  "
$ws.Range("D16").Value = "'13.53"
$ws.Range("E16").Value = "  -8.85%  "
$ws.Range("D17").Value = "'5.35"
$ws.Range("E17").Value = "  -5.75%  "
$ws.Range("D18").Value = "1.964.91"
$ws.Range("E18").Value = "  -4.55%  "
$ws.Range("D19").Value = "36.290.27"
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("D20").Value = "'71.44"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("D22").Value = "'5.13"
$ws.Range("E22").Value = "  -6.62%  "
$ws.Range("D23").Value = "'231.01"
$ws.Range("E23").Value = "  -3.79%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("D26").Value = "'2.27"
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("D27").Value = "'9.59"
$ws.Range("E27").Value = "  -4.78%  "
$ws.Range("D28").Value = "'165.55"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").Value = "'19.89"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("D30").Value = "'0.122"
$ws.Range("E30").Value = "  -4.75%  "
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("E33").Value = "  -8.50%  "
$ws.Range("D34").Value = "'0.0640"
$ws.Range("E34").Value = "  +1.54%  "
$ws.Range("D35").Value = "'4.32"
$ws.Range("E35").Value = "  -8.04%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'1.82"
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("E38").Value = "  -7.52%  "
$ws.Range("E39").Value = "  -13.17%  "
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("D41").Value = "'0.0962"
$ws.Range("E41").Value = "  -6.46%  "
$ws.Range("E42").Value = "  -4.80%  "
$ws.Range("E43").Value = "  -8.50%  "
$ws.Range("E44").Value = "  -4.15%  "
$ws.Range("E45").Value = "  -9.67%  "
$ws.Range("D46").Value = "'15.74"
$ws.Range("E46").Value = "  -9.23%  "
$ws.Range("D47").Value = "'88.88"
$ws.Range("E47").Value = "  -7.12%  "
$ws.Range("D48").Value = "1.348.42"
$ws.Range("E48").Value = "  -3.82%  "
$ws.Range("D49").Value = "'7.28"
$ws.Range("E49").Value = "  -8.71%  "
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("D51").Value = "'44.86"
$ws.Range("E51").Value = "  -4.01%  "

# Reset number-format/quote-prefix styling back to the original default style
# for cells where a leading apostrophe was needed to preserve exact text.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"

Write-Output "Applied cryptos list update"
